# "added item sums, updated UI"
#
# 1. Rename the existing sheet "Sheet2" -> "Skill Sums"
# 2. Add a brand new sheet "Item Sums" right after it, make it the active tab
# 3. Populate "Item Sums" with the Health_Pot / Mana_Pot / Used summary table
# 4. Leave a selection on the new sheet at E8 (matches authored workbook)

$wb = $excel.ActiveWorkbook

# --- 1. rename the original (only) worksheet ---------------------------
$skillSums = $wb.Worksheets.Item(1)
$skillSums.Name = "Skill Sums"

# --- 2. insert the new worksheet after "Skill Sums" ---------------------
$itemSums = $wb.Worksheets.Add($null, $skillSums)
$itemSums.Name = "Item Sums"

# --- 3. fill in the data, column by column so shared-strings are --------
#        interned in the same order as the authored workbook
# Column A - item type labels
$itemSums.Range("A1").Value = "Health_Pot"
$itemSums.Range("A2").Value = "Health_Pot"
$itemSums.Range("A3").Value = "Health_Pot"
$itemSums.Range("A4").Value = "Mana_Pot"
$itemSums.Range("A5").Value = "Mana_Pot"
$itemSums.Range("A6").Value = "Mana_Pot"
$itemSums.Range("A7").Value = "Used "

# Column B - item counter values
$itemSums.Range("B1").Value = 552132
$itemSums.Range("B2").Value = 574965
$itemSums.Range("B3").Value = 575066
$itemSums.Range("B4").Value = 568725
$itemSums.Range("B5").Value = 591551
$itemSums.Range("B6").Value = 591647
$itemSums.Range("B7").Value = "700000+"

# Column D - conclusion / legend notes
$itemSums.Range("D1").Value = "Conclusion:"
$itemSums.Range("D2").Value = "Health_Pot = 0-559000 & 570000-589000"
$itemSums.Range("D3").Value = "Mana_Pot = 560000 & 590000-690000"
$itemSums.Range("D4").Value = "Used = 700000+"

# --- 4. match the authored selection / active cell on the new sheet -----
$itemSums.Range("E8").Select() | Out-Null
